# Update the "Trans Tables" sheet: change the E8:E37 formula so that it
# divides by the average of J8:J37 (absolute reference) instead of the
# row-relative J cell, and update the active selection to L15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trans Tables")
$ws.Activate()

$ws.Range("E8").Formula = '="*"&1/30/AVERAGE($J$8:$J$37)'
$ws.Range("E9:E37").Formula = '="*"&1/30/AVERAGE($J$8:$J$37)'

$ws.Range("L15").Select()
